$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.09938814878539404
$ws.Range("D2").Value = 0.1384643024135066
$ws.Range("E2").Value = 0.1155545705233969
$ws.Range("F2").Value = 2.542000168866508
$ws.Range("G2").Value = 0.002561627388379081
$ws.Range("I2").Value = 2.367744596680183
$ws.Range("K2").Value = 2.750437281852953
$ws.Range("L2").Value = 0.1946243542569732

$ws.Range("C3").Value = 0.09722501049695609
$ws.Range("D3").Value = 0.1401662198227491
$ws.Range("E3").Value = 0.1129609120015225
$ws.Range("F3").Value = 2.516803423652647
$ws.Range("G3").Value = 0.002567701942429623
$ws.Range("I3").Value = 2.349219118362356
$ws.Range("K3").Value = 2.579298125689718
$ws.Range("L3").Value = 0.1901722168959665

$ws.Range("C4").Value = 0.09593890276620698
$ws.Range("D4").Value = 0.1412723616166431
$ws.Range("E4").Value = 0.1114343895490926
$ws.Range("F4").Value = 2.503081329172005
$ws.Range("G4").Value = 0.002571623484156849
$ws.Range("I4").Value = 2.339411142944215
$ws.Range("K4").Value = 2.475668663549243
$ws.Range("L4").Value = 0.1875624879319773

$ws.Range("C5").Value = 0.09542541908578528
$ws.Range("D5").Value = 0.1417384400757804
$ws.Range("E5").Value = 0.110828854844069
$ws.Range("F5").Value = 2.497926094726637
$ws.Range("G5").Value = 0.002573269936149855
$ws.Range("I5").Value = 2.335805700921526
$ws.Range("K5").Value = 2.433801994056921
$ws.Range("L5").Value = 0.186529969231195

$ws.Range("C6").Value = 0.09534079780871707
$ws.Range("D6").Value = 0.1418167553378602
$ws.Range("E6").Value = 0.1107293032863552
$ws.Range("F6").Value = 2.497096352436643
$ws.Range("G6").Value = 0.00257354625610795
$ws.Range("I6").Value = 2.335230583653185
$ws.Range("K6").Value = 2.426871955757804
$ws.Range("L6").Value = 0.1863603851596096

$ws.Range("C7").Value = 0.09593193471174288
$ws.Range("D7").Value = 0.1412785853783447
$ws.Range("E7").Value = 0.1114261562207588
$ws.Range("F7").Value = 2.503010040077271
$ws.Range("G7").Value = 0.00257164549231672
$ws.Range("I7").Value = 2.33936093736871
$ws.Range("K7").Value = 2.475102566549481
$ws.Range("L7").Value = 0.1875484378860222

$ws.Range("C8").Value = 0.09863358842538616
$ws.Range("D8").Value = 0.1390383930877412
$ws.Range("E8").Value = 0.1146465426382868
$ws.Range("F8").Value = 2.532947468106542
$ws.Range("G8").Value = 0.002563682205826404
$ws.Range("I8").Value = 2.361030263834792
$ws.Range("K8").Value = 2.691126359242332
$ws.Range("L8").Value = 0.193063425137538

$ws.Range("C9").Value = 0.1042641407325817
$ws.Range("D9").Value = 0.1351333662681675
$ws.Range("E9").Value = 0.1214885707981708
$ws.Range("F9").Value = 2.605677939921861
$ws.Range("G9").Value = 0.002549579492072017
$ws.Range("I9").Value = 2.41607705891947
$ws.Range("K9").Value = 3.126360612563076
$ws.Range("L9").Value = 0.2048708262859265

$ws.Range("C10").Value = 0.1086027980894784
$ws.Range("D10").Value = 0.1325655115336843
$ws.Range("E10").Value = 0.1268418636868063
$ws.Range("F10").Value = 2.667875240460262
$ws.Range("G10").Value = 0.002540129362651458
$ws.Range("I10").Value = 2.464349452734439
$ws.Range("K10").Value = 3.45338991373211
$ws.Range("L10").Value = 0.2141651573040377

$ws.Range("C11").Value = 0.110620296121823
$ws.Range("D11").Value = 0.1314635281394771
$ws.Range("E11").Value = 0.1293493226709401
$ws.Range("F11").Value = 2.698120026806919
$ws.Range("G11").Value = 0.002536025656196583
$ws.Range("I11").Value = 2.488048886502128
$ws.Range("K11").Value = 3.603784991809846
$ws.Range("L11").Value = 0.2185311290787837

$ws.Range("C12").Value = 0.1113905535190582
$ws.Range("D12").Value = 0.1310558185355042
$ws.Range("E12").Value = 0.1303093057392815
$ws.Range("F12").Value = 2.709857181509392
$ws.Range("G12").Value = 0.002534499571942409
$ws.Range("I12").Value = 2.497276540349532
$ws.Range("K12").Value = 3.660972732154335
$ws.Range("L12").Value = 0.2202044892143391

$ws.Range("C13").Value = 0.1112243861415578
$ws.Range("D13").Value = 0.131143198464386
$ws.Range("E13").Value = 0.1301020896093306
$ws.Range("F13").Value = 2.707316682051186
$ws.Range("G13").Value = 0.002534827003495302
$ws.Range("I13").Value = 2.495277891477613
$ws.Range("K13").Value = 3.648645777468403
$ws.Range("L13").Value = 0.2198432052758506

$ws.Range("C14").Value = 0.1106835400526762
$ws.Range("D14").Value = 0.1314297930376078
$ws.Range("E14").Value = 0.129428090911567
$ws.Range("F14").Value = 2.699079934997059
$ws.Range("G14").Value = 0.002535899546263632
$ws.Range("I14").Value = 2.488802959636331
$ws.Range("K14").Value = 3.608485110364484
$ws.Range("L14").Value = 0.2186683940491321

$ws.Range("C15").Value = 0.1103530725993664
$ws.Range("D15").Value = 0.131606591346074
$ws.Range("E15").Value = 0.1290166121613581
$ws.Range("F15").Value = 2.694071795820662
$ws.Range("G15").Value = 0.002536560138021587
$ws.Range("I15").Value = 2.484869938862786
$ws.Range("K15").Value = 3.583916397629139
$ws.Range("L15").Value = 0.2179514076423601

$ws.Range("C16").Value = 0.1084718293042499
$ws.Range("D16").Value = 0.1326388662309839
$ws.Range("E16").Value = 0.1266794551660269
$ws.Range("F16").Value = 2.6659382343625
$ws.Range("G16").Value = 0.002540401462912376
$ws.Range("I16").Value = 2.462835907043853
$ws.Range("K16").Value = 3.443594192075579
$ws.Range("L16").Value = 0.2138826256330049

$ws.Range("C17").Value = 0.1073289513770845
$ws.Range("D17").Value = 0.1332891287263926
$ws.Range("E17").Value = 0.125264236608416
$ws.Range("F17").Value = 2.649181411252187
$ws.Range("G17").Value = 0.002542807867170414
$ws.Range("I17").Value = 2.44976652138547
$ws.Range("K17").Value = 3.357929589799483
$ws.Range("L17").Value = 0.2114220476512543

$ws.Range("C18").Value = 0.1066757243243899
$ws.Range("D18").Value = 0.133669363946666
$ws.Range("E18").Value = 0.1244570331304686
$ws.Range("F18").Value = 2.639726744024131
$ws.Range("G18").Value = 0.00254421035098068
$ws.Range("I18").Value = 2.442412961252629
$ws.Range("K18").Value = 3.308810647408166
$ws.Range("L18").Value = 0.2100197574426659

$ws.Range("C19").Value = 0.1064552623813313
$ws.Range("D19").Value = 0.1337991713055473
$ws.Range("E19").Value = 0.1241848922475555
$ws.Range("F19").Value = 2.63655694956708
$ws.Range("G19").Value = 0.002544688370285899
$ws.Range("I19").Value = 2.439951181442723
$ws.Range("K19").Value = 3.292206039219138
$ws.Range("L19").Value = 0.2095471867097984

$ws.Range("C20").Value = 0.1074501858487906
$ws.Range("D20").Value = 0.1332192627341087
$ws.Range("E20").Value = 0.1254141854555115
$ws.Range("F20").Value = 2.650946195275054
$ws.Range("G20").Value = 0.002542549799807343
$ws.Range("I20").Value = 2.451140825915843
$ws.Range("K20").Value = 3.367032874332949
$ws.Range("L20").Value = 0.2116826364919291

$ws.Range("C21").Value = 0.1108422295682345
$ws.Range("D21").Value = 0.1313453524490882
$ws.Range("E21").Value = 0.1296257760780719
$ws.Range("F21").Value = 2.701491527932802
$ws.Range("G21").Value = 0.002535583758658033
$ws.Range("I21").Value = 2.490697908639376
$ws.Range("K21").Value = 3.620274837246029
$ws.Range("L21").Value = 0.21901291829532

$ws.Range("C22").Value = 0.1130956890220745
$ws.Range("D22").Value = 0.1301765635855929
$ws.Range("E22").Value = 0.1324393170019249
$ws.Range("F22").Value = 2.736183524525387
$ws.Range("G22").Value = 0.002531193593893942
$ws.Range("I22").Value = 2.518027953843458
$ws.Range("K22").Value = 3.78716211386444
$ws.Range("L22").Value = 0.2239207228196562

$ws.Range("C23").Value = 0.1118896380619674
$ws.Range("D23").Value = 0.1307952261796039
$ws.Range("E23").Value = 0.1309320672825578
$ws.Range("F23").Value = 2.717514883016349
$ws.Range("G23").Value = 0.002533521890639674
$ws.Range("I23").Value = 2.503305238915203
$ws.Range("K23").Value = 3.697964257078581
$ws.Range("L23").Value = 0.2212905497532773

$ws.Range("C24").Value = 0.107395363815499
$ws.Range("D24").Value = 0.1332508292494623
$ws.Range("E24").Value = 0.1253463735765052
$ws.Range("F24").Value = 2.650147778975906
$ws.Range("G24").Value = 0.002542666412898896
$ws.Range("I24").Value = 2.450519004062969
$ws.Range("K24").Value = 3.362916872221717
$ws.Range("L24").Value = 0.2115647858924206

$ws.Range("C25").Value = 0.1027054173542155
$ws.Range("D25").Value = 0.1361371608478592
$ws.Range("E25").Value = 0.1195806813906373
$ws.Range("F25").Value = 2.584479215333943
$ws.Range("G25").Value = 0.002553233823545527
$ws.Range("I25").Value = 2.399823372174495
$ws.Range("K25").Value = 3.007356510386387
$ws.Range("L25").Value = 0.2015688854258286
